$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("list")

$data = @(
    @("Apple", 130, "BDT"),
    @("Mango", 43, "BDT"),
    @("Orange", 14, "BDT"),
    @("Banana", 1111, "BDT"),
    @("Pineapple", 48, "BDT")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
